$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Ali Sweets
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Ali Sweets"
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "35645757567"
$ws.Cells.Item(7, 4).Value = 0

# Row 8: KKR
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "KKR"
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "42354354354"
$ws.Cells.Item(8, 4).Value = 0
